# AutoCommit_3 июля 2023 г. 14:41:12_SibNout2020
# Fill previously-blank score cells (C:S, rows 4-31) with explicit 0 values,
# set row 28 (sheet row 31) D:F to 5, shrink that row's height, and move the
# active selection to S1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 (student #1) ---
$ws.Range("D4:S4").Value = 0

# --- Row 5 (student #2) ---
$ws.Range("C5:S5").Value = 0

# --- Row 6 (student #3) ---
$ws.Range("R6:S6").Value = 0

# --- Row 7 (student #4) ---
$ws.Range("E7").Value = 0
$ws.Range("S7").Value = 0

# --- Row 8 (student #5) ---
$ws.Range("D8:S8").Value = 0

# --- Row 9 (student #6) ---
$ws.Range("N9").Value = 0

# --- Row 10 (student #7) ---
$ws.Range("D10:F10").Value = 0
$ws.Range("S10").Value = 0

# --- Row 11 (student #8) ---
$ws.Range("D11:F11").Value = 0
$ws.Range("S11").Value = 0

# --- Row 12 (student #9) ---
$ws.Range("D12:F12").Value = 0

# --- Row 13 (student #10) ---
$ws.Range("R13:S13").Value = 0

# --- Row 14 (student #11) ---
$ws.Range("D14:S14").Value = 0

# --- Row 15 (student #12) ---
$ws.Range("R15:S15").Value = 0

# --- Row 16 (student #13) ---
$ws.Range("R16:S16").Value = 0

# --- Row 17 (student #14) ---
$ws.Range("R17:S17").Value = 0

# --- Row 18 (student #15) ---
$ws.Range("D18:F18").Value = 0
$ws.Range("S18").Value = 0

# --- Row 20 (student #17) ---
$ws.Range("R20:S20").Value = 0

# --- Row 23 (student #20) ---
$ws.Range("R23:S23").Value = 0

# --- Row 25 (student #22) ---
$ws.Range("F25").Value = 0
$ws.Range("P25").Value = 0

# --- Row 27 (student #24) ---
$ws.Range("S27").Value = 0

# --- Row 29 (student #26) ---
$ws.Range("G29").Value = 0
$ws.Range("R29").Value = 0

# --- Row 30 (student #27) ---
$ws.Range("C30:F30").Value = 0
$ws.Range("H30").Value = 0
$ws.Range("K30:L30").Value = 0

# --- Row 31 (student #28) ---
$ws.Range("D31:F31").Value = 5
$ws.Range("R31:S31").Value = 0

# Shrink row 31's height (it no longer needs to wrap as much text).
$ws.Rows.Item(31).RowHeight = 22

# Move the active selection to S1.
$ws.Range("S1").Select()
